# Updated cryptos list on Tue Dec  5 04:45:55 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and fixes the swapped TerraClassic/VeChain rows (41/42).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text still "looks like" a plain number (e.g. 231.99) are
# forced to a text format before/after the write so Excel keeps them as
# literal strings (matching the workbook's inlineStr text cells) instead of
# silently re-typing them as numeric values.
$ws.Range("D2").Value = '41.790.21'
$ws.Range("E2").Value = '  +2.27%  '
$ws.Range("D3").Value = '2.228.42'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.45%  '
$ws.Range("E6").Value = '  -1.64%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.58'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -7.04%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.405'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.28'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("E12").Value = '  -0.27%  '
$ws.Range("D13").Value = '2.560.40'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.39%  '
$ws.Range("E16").Value = '  -3.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("D18").Value = '2.242.36'
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").Value = '41.732.74'
$ws.Range("E19").Value = '  +2.43%  '
$ws.Range("E20").Value = '  +0.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.51'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.13%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.44%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.93'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.27%  '
$ws.Range("E24").Value = '  -0.15%  '
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("E26").Value = '  +0.39%  '
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.28'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.23%  '
$ws.Range("E29").Value = '  -3.00%  '
$ws.Range("E30").Value = '  -2.20%  '
$ws.Range("E31").Value = '  -2.97%  '
$ws.Range("E32").Value = '  -8.34%  '
$ws.Range("E33").Value = '  -1.80%  '
$ws.Range("E34").Value = '  +3.61%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0653'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.04%  '
$ws.Range("E37").Value = '  -8.98%  '
$ws.Range("E38").Value = '  -5.41%  '
$ws.Range("E39").Value = '  -3.83%  '
$ws.Range("E40").Value = '  +0.08%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0240'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.74%  '
$ws.Range("B42").Value = 'TerraClassic'
$ws.Range("C42").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.000235'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.49%  '
$ws.Range("E43").Value = '  -1.79%  '
$ws.Range("E44").Value = '  -1.77%  '
$ws.Range("E45").Value = '  -3.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.49%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0963'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("D48").Value = '1.470.44'
$ws.Range("E48").Value = '  -3.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.48%  '
$ws.Range("E50").Value = '  +9.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.19%  '
